$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 4
$ws.Range("H9").Value = 4
$ws.Range("E10").Value = 207
$ws.Range("F10").Value = 89
$ws.Range("H10").Value = 89
$ws.Range("E11").Value = 154
$ws.Range("F11").Value = 77
$ws.Range("H11").Value = 77
$ws.Range("E12").Value = 226
$ws.Range("E14").Value = 69
$ws.Range("E16").Value = 87
$ws.Range("F16").Value = 42
$ws.Range("H16").Value = 42
$ws.Range("E17").Value = 43
$ws.Range("E22").Value = 88
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 44
$ws.Range("H23").Value = 44
$ws.Range("E24").Value = 104
$ws.Range("E25").Value = 99
$ws.Range("E26").Value = 61
$ws.Range("F26").Value = 27
$ws.Range("H26").Value = 27
$ws.Range("E27").Value = 150
$ws.Range("F27").Value = 73
$ws.Range("H27").Value = 73
$ws.Range("E28").Value = 92
$ws.Range("E29").Value = 94
$ws.Range("F29").Value = 53
$ws.Range("H29").Value = 53
$ws.Range("E30").Value = 109
$ws.Range("F30").Value = 52
$ws.Range("H30").Value = 52
$ws.Range("E32").Value = 102
$ws.Range("E33").Value = 135
$ws.Range("E34").Value = 105
$ws.Range("E35").Value = 68
$ws.Range("F35").Value = 34
$ws.Range("H35").Value = 34
$ws.Range("E37").Value = 72
$ws.Range("F37").Value = 32
$ws.Range("H37").Value = 32
$ws.Range("E39").Value = 112
$ws.Range("E40").Value = 142
$ws.Range("E41").Value = 184
$ws.Range("F41").Value = 64
$ws.Range("H41").Value = 64
$ws.Range("E42").Value = 156
$ws.Range("E43").Value = 51
$ws.Range("F43").Value = 19
$ws.Range("H43").Value = 19
$ws.Range("E44").Value = 144
$ws.Range("F44").Value = 67
$ws.Range("H44").Value = 67
$ws.Range("E46").Value = 126
$ws.Range("E47").Value = 214
$ws.Range("E48").Value = 105
$ws.Range("F48").Value = 31
$ws.Range("H48").Value = 31
$ws.Range("E51").Value = 102
$ws.Range("F51").Value = 36
$ws.Range("H51").Value = 36
$ws.Range("E52").Value = 7
